# repull data, push all data, mean calculation
#
# Updates the "dSF" column (column F) values for a set of rows to reflect
# re-pulled data. Only column F changes; the other columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = -3
    16 = 4
    21 = 0
    24 = -2
    27 = 1
    28 = -1
    32 = -2
    35 = 2
    50 = 5
    51 = 0
    52 = -4
    53 = -4
    56 = -7
    64 = -2
    65 = -2
    66 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
